# cv125152a.xlsx - "correcao nos dados e inicio da analise PNAD 2009"
#
# The sheet had a spurious header-only row ("grandes regioes e unidades da
# federacao", row 6) that carried no B:G data - the real figures for the
# first data row of that block ("norte") were one row further down. The fix
# is to delete that empty header row so every region's numbers shift up onto
# the correct label, the trailing row (37, "goias") collapses into row 36,
# and the now-unused shared string is dropped automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Delete()
